$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New order row appended below the existing data (row 15 was the prior last row).
# Source values are plain text (SKU/qty/cost columns are all text-typed in this
# sheet), so force Text format before assigning the numeric-looking strings to
# keep them stored as text instead of being auto-coerced into numbers.
$ws.Range("A16:E16").NumberFormat = "@"

$ws.Range("A16").Value = "TN380"
$ws.Range("B16").Value = "Natalie's - Strawberry Lemonade"
$ws.Range("C16").Value = "1"
$ws.Range("D16").Value = "10.15"
$ws.Range("E16").Value = "10.15"

# Match the formatting (default/general style) of the rest of the data rows by
# copying the format from the row above, rather than leaving the explicit
# Text number format applied.
$ws.Range("A15:E15").Copy()
$ws.Range("A16:E16").PasteSpecial(-4122)
$excel.CutCopyMode = $false
